# Updated symbol list on Sat Dec 17 06:50:56 UTC 2022 with GitHub Actions
# Applies the latest price / label refresh to the cryptos sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are stored as text (not numbers) in this sheet, so a
# leading apostrophe is used to force Excel to keep them as text when the
# string looks numeric (matches the original inlineStr formatting).

$ws.Range("D2").Value = "'231.84"
$ws.Range("D3").Value = "'22.74"
$ws.Range("D4").Value = "'5.289"
$ws.Range("D5").Value = "'0.05601"
$ws.Range("D7").Value = "'6.460"
$ws.Range("D8").Value = "'1.067"
$ws.Range("D9").Value = "'0.7832"
$ws.Range("D10").Value = "'0.1379"
$ws.Range("D11").Value = "'0.07397"
$ws.Range("D12").Value = "'0.03152"
$ws.Range("D15").Value = "'0.001656"
$ws.Range("D16").Value = "'3.252"
$ws.Range("D17").Value = "'0.04762"

$ws.Range("D18").Value = "'0.0005789"
$ws.Range("E18").Value = "17OneONE"

$ws.Range("D19").Value = "'0.006242"
$ws.Range("D20").Value = "'0.005239"
$ws.Range("D21").Value = "'0.001055"
$ws.Range("D22").Value = "'0.0001500"
$ws.Range("D23").Value = "'3.974"
$ws.Range("D24").Value = "'2.147"

$ws.Range("D27").Value = "'0.0004999"
$ws.Range("E27").Value = "26UpBotsUBXT"

$ws.Range("D40").Value = "'0.04034"

$ws.Range("D41").Value = "'0.007027"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"

$ws.Range("D43").Value = "'0.003216"
$ws.Range("D44").Value = "'0.009186"
$ws.Range("D45").Value = "'0.00005439"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D47").Value = "'0.7851"

$ws.Range("D48").Value = "'0.04255"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"

$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D50").Value = "'0.01010"
